# Update column C (years 2025-2055, rows 4-34) with new IRR/value figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 3).Value = -294.64969276916855
$ws.Cells.Item(5, 3).Value = 3373.3243207832961
$ws.Cells.Item(6, 3).Value = 4454.1340187283349
$ws.Cells.Item(7, 3).Value = 4252.5544228183471
$ws.Cells.Item(8, 3).Value = 4145.6081918410773
$ws.Cells.Item(9, 3).Value = 3991.2411864889227
$ws.Cells.Item(10, 3).Value = 4647.1365624369237
$ws.Cells.Item(11, 3).Value = 4917.6168475751447
$ws.Cells.Item(12, 3).Value = 4841.2698711477851
$ws.Cells.Item(13, 3).Value = 4827.0636839758909
$ws.Cells.Item(14, 3).Value = 4985.0840320289926
$ws.Cells.Item(15, 3).Value = 4873.3212602541498
$ws.Cells.Item(16, 3).Value = 4864.9725323384228
$ws.Cells.Item(17, 3).Value = 4915.6639471722756
$ws.Cells.Item(18, 3).Value = 4929.0233864648735
$ws.Cells.Item(19, 3).Value = 4838.3388712568885
$ws.Cells.Item(20, 3).Value = 4464.6064739359508
$ws.Cells.Item(21, 3).Value = 7580.6163414785888
$ws.Cells.Item(22, 3).Value = 4214.6664323041505
$ws.Cells.Item(23, 3).Value = 3991.7113527908809
$ws.Cells.Item(24, 3).Value = -5320.1595258045809
$ws.Cells.Item(25, 3).Value = 20314.420745591768
$ws.Cells.Item(26, 3).Value = 1009.6186778125118
$ws.Cells.Item(27, 3).Value = 866.58507918659575
$ws.Cells.Item(28, 3).Value = 880.22212473874549
$ws.Cells.Item(29, 3).Value = 845.99770059757248
$ws.Cells.Item(30, 3).Value = 809.0960804716151
$ws.Cells.Item(31, 3).Value = 1088.1202962717
$ws.Cells.Item(32, 3).Value = 739.26813228975504
$ws.Cells.Item(33, 3).Value = 190.62445762365363
$ws.Cells.Item(34, 3).Value = -3198.0886193465772

# Update the active selection on Sheet1 to S8 (was R3)
$ws.Activate()
$ws.Range("S8").Select()
